# Commit: "I0 and IF added"
# Adds two new columns, I ("I0") and J ("IF"), to the sheet: a header cell
# in row 1 (styled like the existing headers) plus a numeric value for
# every data row (2-70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header row (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-70 for columns I (I0) and J (IF)
$data = @(
    @(2, 9, 9),
    @(3, 7, 7),
    @(4, 8, 8),
    @(5, 9, 9),
    @(6, 9, 9),
    @(7, 9, 9),
    @(8, 9, 9),
    @(9, 10, 10),
    @(10, 9, 9),
    @(11, 9, 9),
    @(12, 9, 9),
    @(13, 7, 7),
    @(14, 9, 9),
    @(15, 9, 9),
    @(16, 9, 9),
    @(17, 9, 10),
    @(18, 9, 9),
    @(19, 7, 7),
    @(20, 9, 9),
    @(21, 8, 8),
    @(22, 10, 10),
    @(23, 9, 9),
    @(24, 7, 7),
    @(25, 9, 9),
    @(26, 9, 9),
    @(27, 9, 9),
    @(28, 9, 9),
    @(29, 8, 9),
    @(30, 7, 7),
    @(31, 8, 8),
    @(32, 9, 9),
    @(33, 9, 9),
    @(34, 9, 9),
    @(35, 8, 9),
    @(36, 7, 7),
    @(37, 9, 9),
    @(38, 9, 9),
    @(39, 9, 9),
    @(40, 9, 9),
    @(41, 8, 9),
    @(42, 9, 9),
    @(43, 9, 9),
    @(44, 9, 9),
    @(45, 7, 7),
    @(46, 9, 9),
    @(47, 10, 10),
    @(48, 9, 9),
    @(49, 9, 9),
    @(50, 8, 8),
    @(51, 8, 9),
    @(52, 9, 9),
    @(53, 9, 9),
    @(54, 9, 9),
    @(55, 9, 9),
    @(56, 9, 9),
    @(57, 9, 9),
    @(58, 9, 10),
    @(59, 9, 9),
    @(60, 7, 7),
    @(61, 9, 9),
    @(62, 7, 7),
    @(63, 9, 9),
    @(64, 9, 9),
    @(65, 7, 7),
    @(66, 6, 6),
    @(67, 8, 8),
    @(68, 5, 5),
    @(69, 3, 3),
    @(70, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
